$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# Column A (NC / student id numbers)
$ws.Range("A2").Value = 20330051920156
$ws.Range("A3").Value = 19330051920133
$ws.Range("A4").Value = 20330051920179

# Column B (Paterno)
$ws.Range("B2").Value = "ANTONIO"
$ws.Range("B3").Value = "NAVARRO"
$ws.Range("B4").Value = "REYES"

# Column C (Materno)
$ws.Range("C2").Value = "GUERRA"
$ws.Range("C3").Value = "HERNANDEZ"
$ws.Range("C4").Value = "SARMIENTO"

# Column D (Nombres)
$ws.Range("D2").Value = "LUIS YAEL"
$ws.Range("D3").Value = "DENISSE MERARY"
$ws.Range("D4").Value = "ESTHER ARISBETH"

# Column E (Nombre_Largo)
$ws.Range("E2").Value = "CONTRIBUYE A LA INTEGRACIÓN Y DESARROLLO DEL PERSONAL EN LA ORGANIZACIÓN"
$ws.Range("E3").Value = "SUPERVISA EL CUMPLIMIENTO DE TAREAS Y PROCESOS PARA EVALUAR LA PRODUCTIVIDAD EN LA ORGANIZACIÓN"
$ws.Range("E4").Value = "CONTRIBUYE A LA INTEGRACIÓN Y DESARROLLO DEL PERSONAL EN LA ORGANIZACIÓN"

# Column F (Grupo)
$ws.Range("F2").Value = "3ARHM"
$ws.Range("F3").Value = "5ARHM"
$ws.Range("F4").Value = "3ARHM"

# Column G (Reprobadas)
$ws.Range("G2").Value = 6
$ws.Range("G3").Value = 6
$ws.Range("G4").Value = 6
